# Actualización automática del mapa (2025-08-27 15:04:57)
# Row 81 (Caso "-566", the duplicate/placeholder entry) was removed from the
# PEBCOM sheet; the subsequent rows (82-84) shift up one position, which is
# exactly what deleting the worksheet row and letting Excel move cells up
# accomplishes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(81).Delete()
